$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the generated student's ID and email (row 2): 123 -> 129
$ws.Range("B2").Value = 129
$ws.Range("C2").Value = "129@qq.com"

# Update the active selection to D4 (as recorded in the saved view state)
$ws.Range("D4").Select()
